$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H7").Value = 177.5342465753425
$ws.Range("I7").Value = 142.027397260274
$ws.Range("J7").Value = 165.6986301369863
$ws.Range("K7").Value = 153.8630136986301
$ws.Range("L7").Value = 165.6986301369863
$ws.Range("M7").Value = 201.2054794520548
$ws.Range("N7").Value = 71.01369863013699
$ws.Range("O7").Value = 319.5616438356165
$ws.Range("G8").Value = 177.5342465753425
$ws.Range("H8").Value = ":"
$ws.Range("I8").Value = ":"
$ws.Range("J8").Value = 35.50684931506849
$ws.Range("K8").Value = 23.67123287671233
$ws.Range("L8").Value = 35.50684931506849
$ws.Range("M8").Value = ":"
$ws.Range("N8").Value = ":"
$ws.Range("O8").Value = 82.84931506849315
$ws.Range("G9").Value = 142.027397260274
$ws.Range("H9").Value = ":"
$ws.Range("I9").Value = ":"
$ws.Range("J9").Value = 35.50684931506849
$ws.Range("K9").Value = 23.67123287671233
$ws.Range("L9").Value = ":"
$ws.Range("M9").Value = 35.50684931506849
$ws.Range("N9").Value = 23.67123287671233
$ws.Range("O9").Value = 23.67123287671233
$ws.Range("G10").Value = 165.6986301369863
$ws.Range("H10").Value = 35.50684931506849
$ws.Range("I10").Value = 35.50684931506849
$ws.Range("J10").Value = ":"
$ws.Range("K10").Value = 11.83561643835616
$ws.Range("L10").Value = 11.83561643835616
$ws.Range("M10").Value = 23.67123287671233
$ws.Range("N10").Value = 11.83561643835616
$ws.Range("O10").Value = 35.50684931506849
$ws.Range("G11").Value = 153.8630136986301
$ws.Range("H11").Value = 23.67123287671233
$ws.Range("I11").Value = 23.67123287671233
$ws.Range("J11").Value = 11.83561643835616
$ws.Range("K11").Value = ":"
$ws.Range("L11").Value = 23.67123287671233
$ws.Range("M11").Value = 23.67123287671233
$ws.Range("N11").Value = 11.83561643835616
$ws.Range("O11").Value = 35.50684931506849
$ws.Range("G12").Value = 165.6986301369863
$ws.Range("H12").Value = 35.50684931506849
$ws.Range("I12").Value = ":"
$ws.Range("J12").Value = 11.83561643835616
$ws.Range("K12").Value = 23.67123287671233
$ws.Range("L12").Value = ":"
$ws.Range("M12").Value = 35.50684931506849
$ws.Range("N12").Value = 11.83561643835616
$ws.Range("O12").Value = 47.34246575342465
$ws.Range("G13").Value = 201.2054794520548
$ws.Range("H13").Value = ":"
$ws.Range("I13").Value = 35.50684931506849
$ws.Range("J13").Value = 23.67123287671233
$ws.Range("K13").Value = 23.67123287671233
$ws.Range("L13").Value = 35.50684931506849
$ws.Range("M13").Value = ":"
$ws.Range("N13").Value = ":"
$ws.Range("O13").Value = 82.84931506849315
$ws.Range("G14").Value = 71.01369863013699
$ws.Range("H14").Value = ":"
$ws.Range("I14").Value = 23.67123287671233
$ws.Range("J14").Value = 11.83561643835616
$ws.Range("K14").Value = 11.83561643835616
$ws.Range("L14").Value = 11.83561643835616
$ws.Range("M14").Value = ":"
$ws.Range("N14").Value = ":"
$ws.Range("O14").Value = 11.83561643835616
$ws.Range("G15").Value = 319.5616438356165
$ws.Range("H15").Value = 82.84931506849315
$ws.Range("I15").Value = 23.67123287671233
$ws.Range("J15").Value = 35.50684931506849
$ws.Range("K15").Value = 35.50684931506849
$ws.Range("L15").Value = 47.34246575342465
$ws.Range("M15").Value = 82.84931506849315
$ws.Range("N15").Value = 11.83561643835616
$ws.Range("O15").Value = ":"
